$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.607064
$ws.Range("H2").Value = 4.821192
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.383078666666667
$ws.Range("N2").Value = 7.149236
$ws.Range("O2").Value = 0.3945674126995297
$ws.Range("P2").Value = 0.4493916190129689
$ws.Range("Q2").Value = 3.829759934368
$ws.Range("R2").Value = 34.467839409312
$ws.Range("S2").Value = 0.3945674126995297
$ws.Range("T2").Value = 0.4493916190129689

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.607064
$ws.Range("H3").Value = 4.821192
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.446170333333334
$ws.Range("N3").Value = 4.338511
$ws.Range("O3").Value = 0.2394430761886234
$ws.Range("P3").Value = 0.2727131238072956
$ws.Range("Q3").Value = 2.324088280568001
$ws.Range("R3").Value = 20.916794525112
$ws.Range("S3").Value = 0.2394430761886234
$ws.Range("T3").Value = 0.2727131238072956

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.607064
$ws.Range("H4").Value = 4.821192
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 2.210476
$ws.Range("N4").Value = 4.420952
$ws.Range("O4").Value = 0.365989511111847
$ws.Range("P4").Value = 0.2778952571797354
$ws.Range("Q4").Value = 3.552376402464
$ws.Range("R4").Value = 21.314258414784
$ws.Range("S4").Value = 0.365989511111847
$ws.Range("T4").Value = 0.2778952571797354
